# Fruta / hortaliza, semanal
# A new daily price record was inserted at row 39 (Macroferia Regional de
# Talca - Haba), pushing the existing rows 39..143 down to 40..144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 39, shifting rows 39-143 -> 40-144
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new record's data
$ws.Cells.Item(39, 1).Value  = 5
$ws.Cells.Item(39, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(39, 3).Value  = "Maule"
$ws.Cells.Item(39, 4).Value  = 45238
$ws.Cells.Item(39, 5).Value  = 7
$ws.Cells.Item(39, 6).Value  = 100112026
$ws.Cells.Item(39, 7).Value  = "Haba"
$ws.Cells.Item(39, 8).Value  = "Sin especificar"
$ws.Cells.Item(39, 9).Value  = "Primera"
$ws.Cells.Item(39, 10).Value = 600
$ws.Cells.Item(39, 11).Value = 8000
$ws.Cells.Item(39, 12).Value = 9000
$ws.Cells.Item(39, 13).Value = 8500
$ws.Cells.Item(39, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(39, 15).Value = "Región del Maule"
$ws.Cells.Item(39, 16).Value = 340
$ws.Cells.Item(39, 17).Value = 25
$ws.Cells.Item(39, 18).Value = "Hortaliza"
